$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2019.3334
$ws.Range("I98").Value = 2026.9048
$ws.Range("K98").Value = 2026.9048
$ws.Range("M98").Value = -528.9048
$ws.Range("H112").Value = 2382.5833
$ws.Range("J112").Value = 2382.5833
$ws.Range("L112").Value = 7147.749899999999
$ws.Range("N112").Value = -9363.749899999999
$ws.Range("H113").Value = 33336184
$ws.Range("I113").Value = 83335690
$ws.Range("J113").Value = 3183.8333
$ws.Range("K113").Value = 83335690
$ws.Range("L113").Value = 3183.8333
$ws.Range("M113").Value = -83332436
$ws.Range("N113").Value = -9691.8333
$ws.Range("H116").Value = 9195.913
$ws.Range("J116").Value = 12713
$ws.Range("L116").Value = 12713
$ws.Range("N116").Value = -19597
$ws.Range("H122").Value = 2019.3334
$ws.Range("I122").Value = 2026.9048
$ws.Range("K122").Value = 6080.7144
$ws.Range("M122").Value = -3630.7144
$ws.Range("H132").Value = 5068.8687
$ws.Range("I132").Value = 5132.086
$ws.Range("K132").Value = 15396.258
$ws.Range("M132").Value = -12866.258
$ws.Range("H137").Value = 2280389.5
$ws.Range("J137").Value = 8363.549999999999
$ws.Range("L137").Value = 25090.65
$ws.Range("N137").Value = -30190.65
$ws.Range("H138").Value = 3811.0896
$ws.Range("I138").Value = 5090.615
$ws.Range("J138").Value = 3503.0557
$ws.Range("K138").Value = 15271.845
$ws.Range("L138").Value = 10509.1671
$ws.Range("M138").Value = -10131.845
$ws.Range("N138").Value = -20789.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6238.3076
$ws.Range("I32").Value = 5554.409
$ws.Range("K32").Value = 5554.409
$ws.Range("M32").Value = -5267.409
$ws.Range("H122").Value = 2077.0625
$ws.Range("I122").Value = 2077.0625
$ws.Range("K122").Value = 6231.1875
$ws.Range("M122").Value = -3781.1875
$ws.Range("H132").Value = 4904392
$ws.Range("I132").Value = 2161.8728
$ws.Range("K132").Value = 6485.6184
$ws.Range("M132").Value = -3955.6184

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2194.1
$ws.Range("I86").Value = 1664
$ws.Range("K86").Value = 1664
$ws.Range("M86").Value = -541
$ws.Range("H89").Value = 2194.1
$ws.Range("I89").Value = 1664
$ws.Range("K89").Value = 8320
$ws.Range("M89").Value = -2704
$ws.Range("H107").Value = 1565.75
$ws.Range("I107").Value = 1479.15
$ws.Range("J107").Value = 1998.75
$ws.Range("K107").Value = 1479.15
$ws.Range("L107").Value = 1998.75
$ws.Range("M107").Value = 440.8499999999999
$ws.Range("N107").Value = -5838.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9015.833000000001
$ws.Range("J99").Value = 3399.3333
$ws.Range("L99").Value = 3399.3333
$ws.Range("N99").Value = -6395.3333
$ws.Range("H122").Value = 2499.3044
$ws.Range("I122").Value = 2201.261
$ws.Range("K122").Value = 6603.782999999999
$ws.Range("M122").Value = -4153.782999999999
$ws.Range("H126").Value = 9015.833000000001
$ws.Range("J126").Value = 3399.3333
$ws.Range("L126").Value = 10197.9999
$ws.Range("N126").Value = -15137.9999
$ws.Range("H132").Value = 10104326
$ws.Range("I132").Value = 2786.4
$ws.Range("J132").Value = 25645154
$ws.Range("K132").Value = 8359.200000000001
$ws.Range("L132").Value = 76935462
$ws.Range("M132").Value = -5829.200000000001
$ws.Range("N132").Value = -76940522
$ws.Range("H134").Value = 3995.6428
$ws.Range("I134").Value = 3994.9167
$ws.Range("K134").Value = 11984.7501
$ws.Range("M134").Value = -9449.750100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 84953.12
$ws.Range("J70").Value = 7421.8887
$ws.Range("L70").Value = 7421.8887
$ws.Range("N70").Value = -7961.8887
$ws.Range("H73").Value = 84953.12
$ws.Range("J73").Value = 7421.8887
$ws.Range("L73").Value = 7421.8887
$ws.Range("N73").Value = -9293.8887
$ws.Range("H107").Value = 5820.7334
$ws.Range("I107").Value = 419.16666
$ws.Range("J107").Value = 9421.777
$ws.Range("K107").Value = 419.16666
$ws.Range("L107").Value = 9421.777
$ws.Range("M107").Value = 1500.83334
$ws.Range("N107").Value = -13261.777
$ws.Range("H113").Value = 1976.2
$ws.Range("I113").Value = 1912.8572
$ws.Range("K113").Value = 1912.8572
$ws.Range("M113").Value = 257.1428000000001
$ws.Range("H122").Value = 3204.093
$ws.Range("I122").Value = 2479.4243
$ws.Range("J122").Value = 5595.5
$ws.Range("K122").Value = 7438.2729
$ws.Range("L122").Value = 16786.5
$ws.Range("M122").Value = -4988.2729
$ws.Range("N122").Value = -21686.5
$ws.Range("H132").Value = 2804.2144
$ws.Range("I132").Value = 1955.7894
$ws.Range("K132").Value = 5867.3682
$ws.Range("M132").Value = -3337.3682

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3748.4
$ws.Range("I7").Value = 3771.8333
$ws.Range("J7").Value = 3713.25
$ws.Range("K7").Value = 3771.8333
$ws.Range("L7").Value = 3713.25
$ws.Range("M7").Value = -3659.8333
$ws.Range("N7").Value = -3937.25
$ws.Range("H40").Value = 4215.15
$ws.Range("I40").Value = 4124.8438
$ws.Range("J40").Value = 4576.375
$ws.Range("K40").Value = 4124.8438
$ws.Range("L40").Value = 4576.375
$ws.Range("M40").Value = -3988.8438
$ws.Range("N40").Value = -4848.375
$ws.Range("H46").Value = 2199.4285
$ws.Range("I46").Value = 2199.5
$ws.Range("K46").Value = 2199.5
$ws.Range("M46").Value = -2011.5
$ws.Range("H61").Value = 16187.75
$ws.Range("I61").Value = 2510.5557
$ws.Range("J61").Value = 33772.715
$ws.Range("K61").Value = 2510.5557
$ws.Range("L61").Value = 33772.715
$ws.Range("M61").Value = -2308.5557
$ws.Range("N61").Value = -34176.715
$ws.Range("H100").Value = 1441.0834
$ws.Range("I100").Value = 1421.6666
$ws.Range("J100").Value = 1499.3334
$ws.Range("K100").Value = 1421.6666
$ws.Range("L100").Value = 1499.3334
$ws.Range("M100").Value = -880.6666
$ws.Range("N100").Value = -2581.3334
$ws.Range("H113").Value = 16187.75
$ws.Range("I113").Value = 2510.5557
$ws.Range("J113").Value = 33772.715
$ws.Range("K113").Value = 2510.5557
$ws.Range("L113").Value = 33772.715
$ws.Range("M113").Value = -340.5556999999999
$ws.Range("N113").Value = -38112.715
$ws.Range("H122").Value = 3863.5908
$ws.Range("I122").Value = 3684.6843
$ws.Range("J122").Value = 4996.6665
$ws.Range("K122").Value = 11054.0529
$ws.Range("L122").Value = 14989.9995
$ws.Range("M122").Value = -8604.052899999999
$ws.Range("N122").Value = -19889.9995
$ws.Range("H126").Value = 3748.4
$ws.Range("I126").Value = 3771.8333
$ws.Range("J126").Value = 3713.25
$ws.Range("K126").Value = 11315.4999
$ws.Range("L126").Value = 11139.75
$ws.Range("M126").Value = -8845.499899999999
$ws.Range("N126").Value = -16079.75
$ws.Range("H132").Value = 4061.525
$ws.Range("I132").Value = 2826.68
$ws.Range("J132").Value = 6119.6
$ws.Range("K132").Value = 8480.039999999999
$ws.Range("L132").Value = 18358.8
$ws.Range("M132").Value = -5950.039999999999
$ws.Range("N132").Value = -23418.8
$ws.Range("H136").Value = 2394.1555
$ws.Range("I136").Value = 2203.6052
$ws.Range("K136").Value = 6610.8156
$ws.Range("M136").Value = -4060.8156

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1930
$ws.Range("I81").Value = 2008.3
$ws.Range("J81").Value = 1773.4
$ws.Range("K81").Value = 4016.6
$ws.Range("L81").Value = 3546.8
$ws.Range("M81").Value = -2955.6
$ws.Range("N81").Value = -5668.8
$ws.Range("H84").Value = 1930
$ws.Range("I84").Value = 2008.3
$ws.Range("J84").Value = 1773.4
$ws.Range("K84").Value = 20083
$ws.Range("L84").Value = 17734
$ws.Range("M84").Value = -14779
$ws.Range("N84").Value = -28342
$ws.Range("H107").Value = 845.53845
$ws.Range("I107").Value = 809.5
$ws.Range("J107").Value = 965.6667
$ws.Range("K107").Value = 2428.5
$ws.Range("L107").Value = 2897.0001
$ws.Range("M107").Value = -508.5
$ws.Range("N107").Value = -6737.0001
$ws.Range("H113").Value = 1385.8889
$ws.Range("I113").Value = 1247.1875
$ws.Range("J113").Value = 1587.6364
$ws.Range("K113").Value = 3741.5625
$ws.Range("L113").Value = 4762.9092
$ws.Range("M113").Value = -1571.5625
$ws.Range("N113").Value = -9102.9092
